$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared text updates (masthead volume/issue number + report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"

# --- Crime-statistics table updates (rows 14-30) ---
$ws.Range("G14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("L14").Value = -57.142857142857
$ws.Range("M14").Value = -72.727272727272
$ws.Range("N14").Value = -93.75
$ws.Range("F15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("N15").Value = -75.510204081632
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -16.666666666666
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 121
$ws.Range("J16").Value = 135
$ws.Range("K16").Value = -10.370370370370
$ws.Range("L16").Value = 45.783132530120
$ws.Range("M16").Value = -9.701492537313
$ws.Range("N16").Value = -90.727969348659
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -37.5
$ws.Range("I17").Value = 256
$ws.Range("J17").Value = 218
$ws.Range("K17").Value = 17.431192660550
$ws.Range("L17").Value = 42.222222222222
$ws.Range("M17").Value = 36.898395721925
$ws.Range("N17").Value = -66.226912928759
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = 7.692307692307
$ws.Range("I18").Value = 100
$ws.Range("J18").Value = 114
$ws.Range("K18").Value = -12.280701754386
$ws.Range("L18").Value = 44.927536231884
$ws.Range("M18").Value = -40.119760479041
$ws.Range("N18").Value = -90.186457311089
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -3.448275862068
$ws.Range("I19").Value = 236
$ws.Range("J19").Value = 248
$ws.Range("K19").Value = -4.838709677419
$ws.Range("L19").Value = 56.291390728476
$ws.Range("M19").Value = 46.583850931677
$ws.Range("N19").Value = -29.552238805970
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -10
$ws.Range("I20").Value = 106
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = 29.268292682926
$ws.Range("L20").Value = 82.758620689655
$ws.Range("M20").Value = 27.710843373494
$ws.Range("N20").Value = -76.496674057649
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -10.344827586206
$ws.Range("F21").Value = 114
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = 9.615384615384
$ws.Range("I21").Value = 834
$ws.Range("J21").Value = 815
$ws.Range("K21").Value = 2.331288343558
$ws.Range("L21").Value = 47.872340425531
$ws.Range("M21").Value = 10.610079575596
$ws.Range("N21").Value = -78.965952080706
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 17
$ws.Range("K22").Value = -22.727272727272
$ws.Range("L22").Value = 70
$ws.Range("M22").Value = 0
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -71.428571428571
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -18.181818181818
$ws.Range("I23").Value = 81
$ws.Range("J23").Value = 74
$ws.Range("K23").Value = 9.459459459459
$ws.Range("L23").Value = 68.75
$ws.Range("M23").Value = 72.340425531914
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 54
$ws.Range("G24").Value = 55
$ws.Range("H24").Value = -1.818181818181
$ws.Range("I24").Value = 421
$ws.Range("J24").Value = 552
$ws.Range("K24").Value = -23.731884057971
$ws.Range("L24").Value = 7.124681933842
$ws.Range("M24").Value = -13.552361396303
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 2.941176470588
$ws.Range("I25").Value = 350
$ws.Range("J25").Value = 308
$ws.Range("K25").Value = 13.636363636363
$ws.Range("L25").Value = 56.25
$ws.Range("M25").Value = -28.571428571428
$ws.Range("C26").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("F26").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 41
$ws.Range("K27").Value = -26.785714285714
$ws.Range("L27").Value = -16.326530612244
$ws.Range("C28").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 14
$ws.Range("K28").Value = -22.222222222222
$ws.Range("L28").Value = -56.25
$ws.Range("M28").Value = -68.181818181818
$ws.Range("N28").Value = -91.566265060241
$ws.Range("C29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -33.333333333333
$ws.Range("I29").Value = 13
$ws.Range("K29").Value = -18.75
$ws.Range("L29").Value = -55.172413793103
$ws.Range("M29").Value = -60.606060606060
$ws.Range("N29").Value = -91.156462585034
$ws.Range("D30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 2
